$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.362.17"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.843.91"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'240.15"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'0.6334"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.07540"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'0.2954"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'24.75"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'0.07729"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "'4.987"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'0.6817"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'82.96"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "'0.000009944"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "'6.132"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "29.388.82"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'230.19"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("D19").Value = "'12.43"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'7.541"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +233.59%  "
$ws.Range("D24").Value = "'16.69"
$ws.Range("E24").Value = "  +172.07%  "
$ws.Range("D25").Value = "'156.41"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'0.1398"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'8.369"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'17.65"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "'2.720"
$ws.Range("E29").Value = "  +172.52%  "
$ws.Range("D30").Value = "'1.469"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "'0.05709"
$ws.Range("D32").Value = "'1.253"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "'4.115"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "'4.015"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'1.839"
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").Value = "'1.154"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").Value = "'0.7139"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "1.241.67"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").Value = "'2.799"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'0.01811"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +264.79%  "
$ws.Range("D43").Value = "'0.9018"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'101.75"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'66.12"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("D47").Value = "'7.054"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").Value = "'9.126"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "'0.4009"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'1.699"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "'0.1120"
$ws.Range("E51").Value = "  -0.92%  "
